$d = $word.ActiveDocument
$d.Content.Find.Execute("Rozk Ai ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Difax ", 2)
